$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the trailing path parameter / query string from each API URL in column B.
$ws.Range("B2").Value = "https://brasilapi.com.br/api/cep/v2"
$ws.Range("B3").Value = "https://brasilapi.com.br/api/ddd/v1"
$ws.Range("B4").Value = "https://brasilapi.com.br/api/banks/v1"
$ws.Range("B5").Value = "https://brasilapi.com.br/api/cnpj/v1"
$ws.Range("B6").Value = "https://brasilapi.com.br/api/ibge/municipios/v1"
$ws.Range("B7").Value = "https://brasilapi.com.br/api/feriados/v1"
$ws.Range("B8").Value = "https://brasilapi.com.br/api/fipe/marcas/v1"
$ws.Range("B9").Value = "https://brasilapi.com.br/api/isbn/v1"
$ws.Range("B10").Value = "https://brasilapi.com.br/api/registrobr/v1"
$ws.Range("B11").Value = "https://brasilapi.com.br/api/taxas/v1"

$ws.Range("B9").Select()
